# Fruta / hortaliza, semanal
# Update the weekly price rows: swap the date/volume/price figures between
# rows 2-3 and rows 4-5 so that the most recent week (44574) moves to the
# top and the earlier week (44223) follows it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was 44574 / 200 / 6000 / 7000 / 6500 / 3250) -> becomes 44223 / 100 / 3500 / 4000 / 3750 / 1875
$ws.Range("D2").Value = 44223
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 3500
$ws.Range("O2").Value = 4000
$ws.Range("P2").Value = 3750
$ws.Range("S2").Value = 1875

# Row 3 (was 44574 / 100 / 5000 / 5000 / 5000 / 2500) -> becomes 44223 / 50 / 3000 / 3000 / 3000 / 1500
$ws.Range("D3").Value = 44223
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("S3").Value = 1500

# Row 4 (was 44223 / 100 / 3500 / 4000 / 3750 / 1875) -> becomes 44574 / 200 / 6000 / 7000 / 6500 / 3250
$ws.Range("D4").Value = 44574
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6500
$ws.Range("S4").Value = 3250

# Row 5 (was 44223 / 50 / 3000 / 3000 / 3000 / 1500) -> becomes 44574 / 100 / 5000 / 5000 / 5000 / 2500
$ws.Range("D5").Value = 44574
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("S5").Value = 2500
